# update concise_ms csv pattern
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): right-answer mark 4 -> 5, wrong-answer penalty -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 ("Total"): total marks 112 -> 140, wrong total 0 -> -0, score label updated
$ws.Range("B12").Value = 140
$ws.Range("C12").Value = -0
$ws.Range("E12").Value = "140.0/140"
